$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a refresh date (serial 45205 = 2023-10-06)
# that was bumped by one day (45206 = 2023-10-07) for every data row (2-387).
$ws.Range("C2:C387").Value = 45206
